# Applies the "corrección de estilo" edits described by the diff to
# GuiaDidactica_CN_06_12_CO.docx
#
# Strategy: use Find/Replace (Range.Find.Execute) against the whole
# document content for each textual change. Find/Replace in Word
# operates on the logical text stream and is agnostic to run
# boundaries, so it correctly matches phrases that are split across
# multiple <w:r> runs in the underlying XML (e.g. "podrán definir" /
# "qu" / "é" / " " / "es la energía...").
#
# wdReplace = 2 (wdReplaceAll), per the Execute(..., Replace) parameter.

$d = $word.ActiveDocument
$nbsp = [char]0x00A0

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute(
        $find,    # FindText
        $true,    # MatchCase
        $false,   # MatchWholeWord
        $false,   # MatchWildcards
        $false,   # MatchSoundsLike
        $false,   # MatchAllWordForms
        $true,    # Forward
        1,        # Wrap (wdFindContinue)
        $false,   # Format
        $replace, # ReplaceWith
        2         # Replace (wdReplaceAll)
    )
    if (-not $ok) {
        Write-Host "NOT FOUND:" $find
    }
}

# 1. Estandar paragraph: add comma before "y explico"
Replace-Text `
    "Verifico la acción de fuerzas electrostáticas y magnéticas y explico su relación con la carga eléctrica." `
    "Verifico la acción de fuerzas electrostáticas y magnéticas, y explico su relación con la carga eléctrica."

# 2. "que es" -> "qué es"; comma -> semicolon before "asimismo"; comma added after "asimismo"
Replace-Text `
    "podrán definir que es la energía y enumerar los tipos existentes, asimismo " `
    "podrán definir qué es la energía y enumerar los tipos existentes; asimismo, "

# 3. Add a trailing space run after "...estos recursos."
Replace-Text `
    "contaminación por el uso indebido de estos recursos." `
    "contaminación por el uso indebido de estos recursos. "

# 4. "se le ofrecen le ayudan" -> "se ofrecen le ayudarán"; also normalizes
#    two stray non-breaking spaces to regular spaces.
Replace-Text `
    ("Los recursos que se le ofrecen le ayudan a" + $nbsp + "profundizar en el conocimiento de los tipos de energía existentes en" + $nbsp + "el entorno inmediato para tomar conciencia ") `
    "Los recursos que se ofrecen le ayudarán a profundizar en el conocimiento de los tipos de energía existentes en el entorno inmediato para tomar conciencia "

# 5. "relacionar las" -> "relación de las"
Replace-Text `
    ", relacionar las " `
    ", relación de las "

# 6. Normalize stray NBSPs; add comma after "importantes"; "estas fuentes de
#    energía" -> "estas"; "medio ambiente" -> "medioambiente"
Replace-Text `
    ("Para" + $nbsp + "que el alumno trabaje" + $nbsp + "en grupo o individualmente, se" + $nbsp + "le presentan una serie de actividades que abarcan otros conceptos importantes como diferenciar los tipos de energía, construir una tabla resumen de las principales fuentes de energía y concientizarse de que estas fuentes de energía pueden contaminar el medio ambiente.") `
    "Para que el alumno trabaje en grupo o individualmente, se le presentan una serie de actividades que abarcan otros conceptos importantes, como diferenciar los tipos de energía, construir una tabla resumen de las principales fuentes de energía y concientizarse de que estas pueden contaminar el medioambiente."

# 7. comma -> semicolon before "la adquisición de"
Replace-Text `
    ", la adquisición de " `
    "; la adquisición de "

# 8. Rework punctuation of the closing sentence of that paragraph; also
#    normalizes a stray non-breaking space and adds "el" before "respeto"
Replace-Text `
    (", la habilidad para buscar, obtener y procesar información en distintos soportes, comprender el valor de las cosas" + $nbsp + "y la tolerancia y respeto por los demás.") `
    "; la habilidad para buscar, obtener y procesar información en distintos soportes; comprender el valor de las cosas, y la tolerancia y el respeto por los demás."

# 9. "potencilizan" -> "potencializa"
Replace-Text `
    " potencilizan la " `
    " potencializa la "

# 10. "Los debates sobre" -> "Los debates acerca de"
Replace-Text `
    ". Los debates sobre las normas de seguridad de las diferentes centrales, así como sobre el consumo responsable de la energía refuerzan la " `
    ". Los debates acerca de las normas de seguridad de las diferentes centrales, así como sobre el consumo responsable de la energía refuerzan la "

# 11. Remove redundant "le" before "ofrecen"; normalize a stray NBSP
Replace-Text `
    ("Las propuestas conceptuales" + $nbsp + "y de carácter práctico le ofrecen la posibilidad de adaptar el discurso en función de las características del grupo. Podrá escoger entre los distintos recursos y actividades con el fin de atender mejor a la ") `
    "Las propuestas conceptuales y de carácter práctico ofrecen la posibilidad de adaptar el discurso en función de las características del grupo. Podrá escoger entre los distintos recursos y actividades con el fin de atender mejor a la "

Write-Host "Done."
